$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1 = 45
    2 = 52
    3 = 91
    4 = 104
    5 = 126
    6 = 155
    7 = 180
    8 = 254
    9 = 288
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
